# Scheduled runner update: refresh Universalis market price / profit figures
$wb = $excel.ActiveWorkbook

# ALC!row 40 - "Stuck in the Moment"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1603.6471
$ws.Range("I40").Value = 1400
$ws.Range("J40").Value = 1784.6666
$ws.Range("K40").Value = 1400
$ws.Range("L40").Value = 1784.6666
$ws.Range("M40").Value = -1225
$ws.Range("N40").Value = -2134.6666

# ALC!row 64 - "Forged from the Void"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3738.9246
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 3723.26
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 3723.26
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -4219.26

# ALC!row 67 - "Dodging the Draft (L)"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3738.9246
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 3723.26
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 3723.26
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -5439.26

# ALC!row 76 - "Warding Off Temptation"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4613.25
$ws.Range("I76").Value = 3696.75
$ws.Range("J76").Value = 5988
$ws.Range("K76").Value = 3696.75
$ws.Range("L76").Value = 5988
$ws.Range("M76").Value = -3381.75
$ws.Range("N76").Value = -6618

# ALC!row 79 - "The Garden of Arcane Delights (L)"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4613.25
$ws.Range("I79").Value = 3696.75
$ws.Range("J79").Value = 5988
$ws.Range("K79").Value = 3696.75
$ws.Range("L79").Value = 5988
$ws.Range("M79").Value = -2604.75
$ws.Range("N79").Value = -8172

# ARM!row 63 - "Rivets Run through It"
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1451559.8
$ws.Range("I63").Value = 1853560.1
$ws.Range("J63").Value = 4358.4
$ws.Range("K63").Value = 1853560.1
$ws.Range("L63").Value = 4358.4
$ws.Range("M63").Value = -1852874.1
$ws.Range("N63").Value = -5730.4

# ARM!row 66 - "A Riveting Revival (L)"
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1451559.8
$ws.Range("I66").Value = 1853560.1
$ws.Range("J66").Value = 4358.4
$ws.Range("K66").Value = 9267800.5
$ws.Range("L66").Value = 21792
$ws.Range("M66").Value = -9264368.5
$ws.Range("N66").Value = -28656

# ARM!row 88 - "The Mast Chance"
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 8266.954
$ws.Range("I88").Value = 1650.2
$ws.Range("J88").Value = 10213.059
$ws.Range("K88").Value = 1650.2
$ws.Range("L88").Value = 10213.059
$ws.Range("M88").Value = -1244.2
$ws.Range("N88").Value = -11025.059

# ARM!row 91 - "The Rose and the Riveter (L)"
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 8266.954
$ws.Range("I91").Value = 1650.2
$ws.Range("J91").Value = 10213.059
$ws.Range("K91").Value = 1650.2
$ws.Range("L91").Value = 10213.059
$ws.Range("M91").Value = -246.2
$ws.Range("N91").Value = -13021.059

# BSM!row 105 - "Ingot to Wing It"
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2146.7576
$ws.Range("I105").Value = 2015.5
$ws.Range("K105").Value = 2015.5
$ws.Range("M105").Value = -268.5

# CRP!row 62 - "Splinter in the Sewers"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5100.385
$ws.Range("I62").Value = 5170.5
$ws.Range("J62").Value = 4866.6665
$ws.Range("K62").Value = 5170.5
$ws.Range("L62").Value = 4866.6665
$ws.Range("M62").Value = -4546.5
$ws.Range("N62").Value = -6114.6665

# CRP!row 65 - "The Lumber of Their Discontent (L)"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5100.385
$ws.Range("I65").Value = 5170.5
$ws.Range("J65").Value = 4866.6665
$ws.Range("K65").Value = 25852.5
$ws.Range("L65").Value = 24333.3325
$ws.Range("M65").Value = -22732.5
$ws.Range("N65").Value = -30573.3325

# GSM!row 70 - "Sky Is the Limit"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 47590.523
$ws.Range("I70").Value = 128536
$ws.Range("J70").Value = 4419.6
$ws.Range("K70").Value = 128536
$ws.Range("L70").Value = 4419.6
$ws.Range("M70").Value = -128266
$ws.Range("N70").Value = -4959.6

# GSM!row 73 - "Hulls of Broken Dreams (L)"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 47590.523
$ws.Range("I73").Value = 128536
$ws.Range("J73").Value = 4419.6
$ws.Range("K73").Value = 128536
$ws.Range("L73").Value = 4419.6
$ws.Range("M73").Value = -127600
$ws.Range("N73").Value = -6291.6

# GSM!row 80 - "Needs More Prayerbell"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2584.2856
$ws.Range("I80").Value = 2348.3333
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2348.3333
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1350.3333
$ws.Range("N80").Value = -5996

# GSM!row 83 - "With a Noise That Reaches Heaven (L)"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2584.2856
$ws.Range("I83").Value = 2348.3333
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 11741.6665
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -6749.666499999999
$ws.Range("N83").Value = -29984

# LTW!row 16 - "Saddle Sore"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6324.8335
$ws.Range("I16").Value = 1544.7
$ws.Range("J16").Value = 30225.5
$ws.Range("K16").Value = 1544.7
$ws.Range("L16").Value = 30225.5
$ws.Range("M16").Value = -1374.7
$ws.Range("N16").Value = -30565.5

# LTW!row 22 - "Skin off Their Backs"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1524.8846
$ws.Range("I22").Value = 1820.0667
$ws.Range("J22").Value = 1122.3636
$ws.Range("K22").Value = 1820.0667
$ws.Range("L22").Value = 1122.3636
$ws.Range("M22").Value = -1525.0667
$ws.Range("N22").Value = -1712.3636

# LTW!row 27 - "Fire and Hide"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1524.8846
$ws.Range("I27").Value = 1820.0667
$ws.Range("J27").Value = 1122.3636
$ws.Range("K27").Value = 1820.0667
$ws.Range("L27").Value = 1122.3636
$ws.Range("M27").Value = -1713.0667
$ws.Range("N27").Value = -1336.3636

# LTW!row 122 - "Hell on Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 17866272
$ws.Range("I122").Value = 34534.668
$ws.Range("J122").Value = 22729474
$ws.Range("K122").Value = 103604.004
$ws.Range("L122").Value = 68188422
$ws.Range("M122").Value = -101154.004
$ws.Range("N122").Value = -68193322

# WVR!row 62 - "Pride Up in Smoke"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4355.7144
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4355.7144
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4355.7144
$ws.Range("N62").Value = -5603.7144
$ws.Range("M62").ClearContents()

# WVR!row 65 - "Desperate for Diversionaries (L)"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4355.7144
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4355.7144
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 21778.572
$ws.Range("N65").Value = -28018.572
$ws.Range("M65").ClearContents()

# WVR!row 81 - "Where the Dragonflies, the Net Catches"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4168138.8
$ws.Range("J81").Value = 2096
$ws.Range("L81").Value = 4192
$ws.Range("N81").Value = -6314

# WVR!row 84 - "To Kill a Dragon on Nameday (L)"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4168138.8
$ws.Range("J84").Value = 2096
$ws.Range("L84").Value = 20960
$ws.Range("N84").Value = -31568

# WVR!row 125 - "Color Coated"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 50490.5
$ws.Range("J125").Value = 50490.5
$ws.Range("L125").Value = 50490.5
$ws.Range("N125").Value = -60330.5
